# Generated PowerShell COM-interop script to transform before.xlsx into the target state.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: A1/B1 lose their thin-border style (C1/D1 keep it) ---
$ws.Range("A1:B1").ClearFormats()
$ws.Range("A1").Value = "ItemID"
$ws.Range("B1").Value = "ItemNameE"
$ws.Range("C1").Value = "NBRUN"
$ws.Range("D1").Value = "PRIXUN"

# --- Rows 40-48 are brand new: clone the bordered style from row 39 first ---
$ws.Range("A39:D39").Copy()
$ws.Range("A40:D48").PasteSpecial(-4122)

# --- Row 38's A cell switches from the left-aligned style (s=2) to the plain bordered style (s=1) ---
$ws.Range("A39").Copy()
$ws.Range("A38").PasteSpecial(-4122)

# --- Data rows 2-48 ---
$ws.Range("A2").Value = 2644573
$ws.Range("B2").Value = 'Bref Triggers Cuisine 500ml'
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = 178
$ws.Range("A3").Value = 2646698
$ws.Range("B3").Value = 'Bref Triggers SDB 500ml'
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = 178
$ws.Range("A4").Value = 2682282
$ws.Range("B4").Value = 'Le Chat Premium RL 2.5L'
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 970
$ws.Range("A5").Value = 2728380
$ws.Range("B5").Value = 'Bref 1.75L Javel desinf'
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 155
$ws.Range("A6").Value = 2728382
$ws.Range("B6").Value = 'Bref 900ml javel disinf'
$ws.Range("C6").Value = 12
$ws.Range("D6").Value = 86.5
$ws.Range("A7").Value = 2766729
$ws.Range("B7").Value = 'ISIS LS gel 900 ml Citron Harmonie  RE                                                                                        '
$ws.Range("C7").Value = 10
$ws.Range("D7").Value = 288
$ws.Range("A8").Value = 2766835
$ws.Range("B8").Value = 'ISIS LS powder 2,5 bag Citron Limitless'
$ws.Range("C8").Value = 4
$ws.Range("D8").Value = 725
$ws.Range("A9").Value = 2791901
$ws.Range("B9").Value = 'Le Chat Regular 1L Adv21 RL'
$ws.Range("C9").Value = 10
$ws.Range("D9").Value = 378
$ws.Range("A10").Value = 2806713
$ws.Range("B10").Value = 'Le Chat Reg 2,5L Adv21 RL'
$ws.Range("C10").Value = 4
$ws.Range("D10").Value = 815
$ws.Range("A11").Value = 2806719
$ws.Range("B11").Value = 'Le Chat Reg 4L Adv21 RL'
$ws.Range("C11").Value = 3
$ws.Range("D11").Value = 1195
$ws.Range("A12").Value = 2817870
$ws.Range("B12").Value = 'Le Chat LS 2,5kg bag Adv21'
$ws.Range("C12").Value = 4
$ws.Range("D12").Value = 795
$ws.Range("A13").Value = 2820870
$ws.Range("B13").Value = 'ISIS HS POWDER ANTIBACTERIAL 750GR LIMITLESS'
$ws.Range("C13").Value = 12
$ws.Range("D13").Value = 203
$ws.Range("A14").Value = 2821686
$ws.Range("B14").Value = 'ISIS HS POWDER ANTIBACTERIAL 300GR LIMITLESS'
$ws.Range("C14").Value = 25
$ws.Range("D14").Value = 86
$ws.Range("A15").Value = 2823409
$ws.Range("B15").Value = 'ISIS HS POWDER SDM 300GR LIMITLESS'
$ws.Range("C15").Value = 25
$ws.Range("D15").Value = 86
$ws.Range("A16").Value = 2823411
$ws.Range("B16").Value = 'ISIS HS POWDER SDM 750 GR LIMITLESS'
$ws.Range("C16").Value = 12
$ws.Range("D16").Value = 203
$ws.Range("A17").Value = 2829475
$ws.Range("B17").Value = 'Le Chat HS 1L RL'
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 331
$ws.Range("A18").Value = 2830113
$ws.Range("B18").Value = 'LE CHAT HS 300 gr'
$ws.Range("C18").Value = 25
$ws.Range("D18").Value = 96
$ws.Range("A19").Value = 2830114
$ws.Range("B19").Value = 'LE CHAT HS 750 gr'
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 225
$ws.Range("A20").Value = 2845955
$ws.Range("B20").Value = 'Pril ISIS Cold Power liquid 1250ml Lemon'
$ws.Range("C20").Value = 12
$ws.Range("D20").Value = 284
$ws.Range("A21").Value = 2845956
$ws.Range("B21").Value = 'Pril ISIS Cold Power liquid 3000ml Lemon'
$ws.Range("C21").Value = 4
$ws.Range("D21").Value = 665
$ws.Range("A22").Value = 2845959
$ws.Range("B22").Value = 'Pril ISIS Cold Power liquid 650ml Lemon'
$ws.Range("C22").Value = 12
$ws.Range("D22").Value = 168
$ws.Range("A23").Value = 2859185
$ws.Range("B23").Value = 'Le Chat Savon de Marseille 2,5L'
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 815
$ws.Range("A24").Value = 2860119
$ws.Range("B24").Value = 'Le Chat Rose LS Gel 2,5L'
$ws.Range("C24").Value = 4
$ws.Range("D24").Value = 815
$ws.Range("A25").Value = 2875891
$ws.Range("B25").Value = 'ISIS HS 300g LEMON LIMITLESS'
$ws.Range("C25").Value = 25
$ws.Range("D25").Value = 86
$ws.Range("A26").Value = 2875892
$ws.Range("B26").Value = 'ISIS HS 750g LEMON LIMITLESS'
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 203
$ws.Range("A27").Value = 2875893
$ws.Range("B27").Value = 'ISIS HS 1,5Kg LEMON LIMITLESS'
$ws.Range("C27").Value = 8
$ws.Range("D27").Value = 397
$ws.Range("A28").Value = 2876884
$ws.Range("B28").Value = 'ISIS LS Gel 2,5L Lemon'
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 495
$ws.Range("A29").Value = 2910278
$ws.Range("B29").Value = 'Nettoyant Moussant Javelisé 900ml'
$ws.Range("C29").Value = 12
$ws.Range("D29").Value = 123.5
$ws.Range("A30").Value = 2917788
$ws.Range("B30").Value = 'Le Chat LS 2,5kg bag Regular'
$ws.Range("C30").Value = 4
$ws.Range("D30").Value = 795
$ws.Range("A31").Value = 2918203
$ws.Range("B31").Value = 'ISIS LS bag 2,5Kg Citron Limitless'
$ws.Range("C31").Value = 4
$ws.Range("D31").Value = 725
$ws.Range("A32").Value = 2922764
$ws.Range("B32").Value = 'Le Chat 2.5L fresco RL'
$ws.Range("C32").Value = 4
$ws.Range("D32").Value = 815
$ws.Range("A33").Value = 2922795
$ws.Range("B33").Value = 'Le Chat 1L fresco RL'
$ws.Range("C33").Value = 10
$ws.Range("D33").Value = 378
$ws.Range("A34").Value = 2940803
$ws.Range("B34").Value = 'Le Chat Rose 1L'
$ws.Range("C34").Value = 10
$ws.Range("D34").Value = 378
$ws.Range("A35").Value = 2940804
$ws.Range("B35").Value = 'Le Chat 4L fresco RL'
$ws.Range("C35").Value = 3
$ws.Range("D35").Value = 1195
$ws.Range("A36").Value = 2943549
$ws.Range("B36").Value = 'Pril ISIS Antibactérien 650ML  Lemon'
$ws.Range("C36").Value = 12
$ws.Range("D36").Value = 168
$ws.Range("A37").Value = 2951994
$ws.Range("B37").Value = 'Le Chat Rose 2,5L ADV23'
$ws.Range("C37").Value = 4
$ws.Range("D37").Value = 1195
$ws.Range("A38").Value = 2952074
$ws.Range("B38").Value = 'Le Chat Regular 2,5L ADV23'
$ws.Range("C38").Value = 4
$ws.Range("D38").Value = 284
$ws.Range("A39").Value = 2952081
$ws.Range("B39").Value = 'Le Chat SDM 2,5L ADV23'
$ws.Range("C39").Value = 4
$ws.Range("D39").Value = 168
$ws.Range("A40").Value = 2952089
$ws.Range("B40").Value = 'Le Chat Regular 1L ADV23'
$ws.Range("C40").Value = 10
$ws.Range("D40").Value = 378
$ws.Range("A41").Value = 2952090
$ws.Range("B41").Value = 'Le Chat Rose 1L ADV23'
$ws.Range("C41").Value = 10
$ws.Range("D41").Value = 378
$ws.Range("A42").Value = 2952095
$ws.Range("B42").Value = 'Le Chat Regular 4L ADV23'
$ws.Range("C42").Value = 3
$ws.Range("D42").Value = 1195
$ws.Range("A43").Value = 2958928
$ws.Range("B43").Value = 'Le Chat 2,5kg bag ADV23'
$ws.Range("C43").Value = 4
$ws.Range("D43").Value = 795
$ws.Range("A44").Value = 2958931
$ws.Range("B44").Value = 'Le Chat HS Gel 1L ADV23'
$ws.Range("C44").Value = 10
$ws.Range("D44").Value = 331
$ws.Range("A45").Value = 2958932
$ws.Range("B45").Value = 'Le Chat HS 300gr ADV23'
$ws.Range("C45").Value = 25
$ws.Range("D45").Value = 96
$ws.Range("A46").Value = 2958933
$ws.Range("B46").Value = 'Le Chat HS 750gr ADV23'
$ws.Range("C46").Value = 12
$ws.Range("D46").Value = 225
$ws.Range("A47").Value = 2970224
$ws.Range("B47").Value = 'Pril Isis ultra power 1.25L'
$ws.Range("C47").Value = 12
$ws.Range("D47").Value = 284
$ws.Range("A48").Value = 2970482
$ws.Range("B48").Value = 'Pril Isis ultra power 650ml'
$ws.Range("C48").Value = 12
$ws.Range("D48").Value = 168

# --- Refresh the view: selection becomes A2:D48, scrolled back to the top ---
$ws.Range("A2:D48").Select()
